$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Personal details ---
Replace-Text "Terpstra" "Morretti"
Replace-Text "Charlotte" "Luca"
Replace-Text "Blijde-In 40, 4473 63 Apeldoorn" "Galleria Vittorio Emanuele II 57, 94327 Reggio Calabria"
Replace-Text "Netherlands" "Italy"
Replace-Text "1995-10-28" "1975-04-06"
Replace-Text "Dutch" "Italian"
Replace-Text "XK8106150" "WR4329993"
Replace-Text "2016-11-10" "2023-10-24"
Replace-Text "2026-11-09" "2033-10-23"

# Gender checkboxes
Replace-Text "☒ Female" "☐ Female"
Replace-Text "☐ Male" "☒ Male"

# Contact details
Replace-Text "06 45709443" "363 9193391"
Replace-Text "charlotte.terpstra@casema.nl" "luca.moretti@tiscali.it"

# Marital status checkboxes
Replace-Text "☒ Single" "☐ Single"
Replace-Text "☐ Widowed" "☒ Widowed"

# Education / employment
Replace-Text "Utrecht University (2019)" "University of Rome La Sapienza (1998)"
Replace-Text "Since 2021" "Since 2007"
Replace-Text "Name Employer Galapagos N.V." "Name Employer Banca Mediolanum S.p.A."
Replace-Text "Position Research Scientist (43000 EUR p.A.)" "Position Wealth Manager (286000 EUR p.A.)"

# Net income checkboxes
Replace-Text "☒ EUR 1.5m-5m " "☐ EUR 1.5m-5m "
Replace-Text "☐ EUR 5m-10m" "☒ EUR 5m-10m"

# Family info
Replace-Text "grandfather,2017,Neurosurgeon" "grandfather,2011,Oil and Gas Executive"

# Income amount
Replace-Text "30000" "480000"

# Net worth checkboxes (< EUR 250,000 -> 250,000 - 500,000)
Replace-Text "☒ < EUR 250,000" "☐ < EUR 250,000"
Replace-Text "☐ EUR 250,000 - 500,000" "☒ EUR 250,000 - 500,000"

# Risk profile checkboxes
Replace-Text "☐ Low" "☒ Low"
Replace-Text "☐ Moderate  ☒ Considerable ☐ High" "☐ Moderate  ☐ Considerable ☐ High"

# Investment horizon checkboxes
Replace-Text "☐ Short" "☒ Short"
Replace-Text "☒ Medium  ☐ Long-Term" "☐ Medium  ☐ Long-Term"

# Asset values
Replace-Text "1660000" "7290000"
Replace-Text "1145400" "4374000"
